$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fasta-method-1")

$ws.Range("C11").Value = 0.0083379745483399999
$ws.Range("D11").Value = 0.0083379745483399999
$ws.Range("E11").Value = 12.99609375
$ws.Range("F11").Value = 143.03515625

$ws.Range("C12").Value = 0.76420497894299999
$ws.Range("D12").Value = 0.772542953491
$ws.Range("E12").Value = 31.125
$ws.Range("F12").Value = 161.03125

$ws.Range("C13").Value = 0.76520609855699995
$ws.Range("D13").Value = 0.77354407310499995
$ws.Range("E13").Value = 31.125
$ws.Range("F13").Value = 161.03125

$ws.Range("C14").Value = 0.00134706497192
$ws.Range("D14").Value = 443.22435212099998
$ws.Range("E14").Value = 239.01953125
$ws.Range("F14").Value = 1101.265625

$ws.Range("C15").Value = 36.746198892599999
$ws.Range("D15").Value = 479.97055101400002
$ws.Range("E15").Value = 239.03515625
$ws.Range("F15").Value = 1101.265625

$ws.Range("C16").Value = 254.481251955
$ws.Range("D16").Value = 734.45180296900003
$ws.Range("E16").Value = 239.046875
$ws.Range("F16").Value = 1101.265625

$ws.Range("C17").Value = 0.00335502624512
$ws.Range("D17").Value = 734.45515799500004
$ws.Range("E17").Value = 239.046875
$ws.Range("F17").Value = 1101.265625

$ws.Range("I25").Select()
